$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45, shifting existing rows 45-72 down to 46-73.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new record.
$ws.Cells.Item(45, 1).Value2 = 10
$ws.Cells.Item(45, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value2 = "La Araucanía"
$ws.Cells.Item(45, 4).Value2 = 44452
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(45, 5).Value2 = 9
$ws.Cells.Item(45, 6).Value2 = 100112031
$ws.Cells.Item(45, 7).Value2 = "Poroto verde"
$ws.Cells.Item(45, 8).Value2 = "Sin especificar"
$ws.Cells.Item(45, 9).Value2 = "Primera"
$ws.Cells.Item(45, 10).Value2 = 50
$ws.Cells.Item(45, 11).Value2 = 40000
$ws.Cells.Item(45, 12).Value2 = 40000
$ws.Cells.Item(45, 13).Value2 = 40000
$ws.Cells.Item(45, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(45, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(45, 16).Value2 = 1600
$ws.Cells.Item(45, 17).Value2 = 25
$ws.Cells.Item(45, 18).Value2 = "Hortaliza"
